$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, shifting existing rows 134:155 down to 135:156
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134
$ws.Cells.Item(134, 1).Value = 5
$ws.Cells.Item(134, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(134, 3).Value = "Maule"
$ws.Cells.Item(134, 4).Value = 45180
$ws.Cells.Item(134, 5).Value = 7
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100107
$ws.Cells.Item(134, 8).Value = "Otros"
$ws.Cells.Item(134, 9).Value = 100107002
$ws.Cells.Item(134, 10).Value = "Chirimoya"
$ws.Cells.Item(134, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(134, 12).Value = "Primera"
$ws.Cells.Item(134, 13).Value = 100
$ws.Cells.Item(134, 14).Value = 23000
$ws.Cells.Item(134, 15).Value = 23000
$ws.Cells.Item(134, 16).Value = 23000
$ws.Cells.Item(134, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(134, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(134, 19).Value = 2300
$ws.Cells.Item(134, 20).Value = 10
